# "forbid delete items and customers when order is in place"
#
# Data cleanup on the orders sheet ("pedidos"):
#   - Rows 49-53 were duplicate placeholder order-lines ("aju" in the
#     Descripcion column). Their customer id (col B) and delivery date
#     (col E) get corrected and the stray placeholder text is cleared.
#   - Rows 54-59 were further duplicate placeholder rows ("aju"/"ajo")
#     that are no longer needed once an order is in place, so they are
#     deleted outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Cliente_id (col B)
$fixups = @{
    49 = 1
    50 = 1
    51 = 1
    52 = 5
    53 = 5
}

foreach ($row in $fixups.Keys) {
    $ws.Cells.Item($row, 2).Value = $fixups[$row]   # Cliente_id
    $ws.Cells.Item($row, 4).ClearContents()          # Descripcion -> blank
    $ws.Cells.Item($row, 5).Value = 45432            # Fecha Entrega
}

# Trailing duplicate placeholder rows are no longer needed; remove them.
$ws.Rows("54:59").Delete()

# Best-effort cleanup of the now-unused custom date-only number formats
# ("yyyy-mm-dd" / "YYYY-MM-DD") that were only applied to the rows just
# removed.
$wb.DeleteNumberFormat("yyyy-mm-dd")
$wb.DeleteNumberFormat("YYYY-MM-DD")
